$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell H1 (bold, bordered,
# centered) onto the two new header cells I1 and J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns: I = I0, J = IF
$iValues = @(6, 8, 8, 5, 6, 7, 1, 11, 5, 1, 5, 4)
$jValues = @(9, 8, 9, 8, 7, 7, 4, 11, 8, 3, 7, 5)

for ($r = 0; $r -lt 12; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
